$d = $word.ActiveDocument

# Locate the "Introducción" paragraph's index.
$introIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd([char]13, [char]7) -eq "Introducción") {
        $introIdx = $i
    }
}

$ns = ' xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$rPr = '<w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri" w:eastAsia="Calibri"/><w:color w:val="auto"/><w:spacing w:val="0"/><w:position w:val="0"/><w:sz w:val="22"/><w:shd w:fill="auto" w:val="clear"/></w:rPr>'
$pPr = '<w:pPr><w:spacing w:before="0" w:after="200" w:line="276"/><w:ind w:right="0" w:left="0" w:firstLine="0"/><w:jc w:val="left"/>' + $rPr + '</w:pPr>'

# Each new line ("-Objetivos", "- Motivación", "- Problema a abordar",
# "- Cómo está organizado") becomes its own paragraph, inserted right after
# "Introducción", in document order.
$runInner = @(
    '<w:t xml:space="preserve"> </w:t><w:tab/><w:t xml:space="preserve">-Objetivos</w:t>',
    '<w:tab/><w:t xml:space="preserve">- Motivación</w:t>',
    '<w:tab/><w:t xml:space="preserve">- Problema a abordar</w:t>',
    '<w:tab/><w:t xml:space="preserve">- Cómo está organizado</w:t>'
)

$afterIdx = $introIdx
foreach ($inner in $runInner) {
    # Insert a fresh paragraph (inherits sibling formatting from the live
    # document) right after the running insertion point.
    $d.Paragraphs($afterIdx).Range.InsertParagraphAfter() | Out-Null
    $afterIdx = $afterIdx + 1

    # Re-fetch the freshly created paragraph's range from the document (the
    # object handed back by InsertParagraphAfter does not reflect further
    # live edits) and inject the run via raw WordprocessingML so a literal
    # <w:tab/> element is produced instead of a tab character in <w:t>.
    $newRange = $d.Paragraphs($afterIdx).Range
    $frag = '<w:p' + $ns + '>' + $pPr + '<w:r>' + $rPr + $inner + '</w:r></w:p>'
    $newRange.InsertXML($frag)

    # InsertXML's WordprocessingML parser normalizes away explicit
    # zero-valued spacing/indentation attributes; restore them explicitly
    # through the paragraph-format object so the saved markup keeps the
    # explicit w:before="0" / w:ind values.
    $pf = $d.Paragraphs($afterIdx).Format
    $pf.SpaceBefore = 0
    $pf.LeftIndent = 0
    $pf.RightIndent = 0
    $pf.FirstLineIndent = 0
}
